$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "OECD" source note block from rows 20-21 down to rows 26-27
# to make room for the new MSME breakdown table (rows 15-19).
$oecdTitle = $ws.Range("A20").Text
$oecdDesc  = $ws.Range("A21").Text
$ws.Range("A20").Clear()
$ws.Range("A21").Clear()

# New table header row (bold "title" style)
$ws.Range("B15").Value = "Number of employees"
$ws.Range("C15").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D15").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B15:D15").Font.Bold = $true

# New table body rows (default/"Normal" style; data cells left blank)
$ws.Range("A16").Value = "Micro"
$ws.Range("A17").Value = "Small"
$ws.Range("A18").Value = "Medium"
$ws.Range("A19").Value = "Large"

# Restore the OECD source note further down the sheet (rows 26-27)
$ws.Range("A26").Value = $oecdTitle
$ws.Range("A26").Font.Bold = $true

$ws.Range("A27").Value = $oecdDesc
$ws.Range("A27").Font.Italic = $true

# Re-assert the pre-existing named-style formatting on the untouched cells
# (the workbook's cell styles are resolved through named cellStyles, so we
# pin the effective font attributes explicitly to keep them intact after
# the save round-trip).
$ws.Range("A1").Font.Size = 18

$ws.Range("A3").Font.Bold = $true
$ws.Range("B9:D9").Font.Bold = $true
$ws.Range("A10").Font.Bold = $true
$ws.Range("A11").Font.Bold = $true

$ws.Range("A7").Font.Bold = $true
$ws.Range("A7").Font.Underline = $true

$ws.Range("A12").Font.Italic = $true
